# Update countries & provincias Spain
#
# This script applies the country-statistics refresh captured in the
# commit: updated case numbers for a handful of countries, two countries
# swap rank order (Georgia/Cabo Verde and Islas Malvinas/Montserrat), and
# the "last updated" timestamp moves from 08:24 to 09:41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 09:41"

# --- India (row 5) ------------------------------------------------------
$ws.Range("B5").Value = 6229474
$ws.Range("C5").Value = 5955
$ws.Range("E5").Value = 944108
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 97541

# --- Rusia (row 7) -------------------------------------------------------
$ws.Range("B7").Value = 1176286
$ws.Range("C7").Value = 8481
$ws.Range("E7").Value = 197307
$ws.Range("G7").Value = 177
$ws.Range("H7").Value = 20722

# --- Israel (row 27) ------------------------------------------------------
$ws.Range("B27").Value = 239222
$ws.Range("C27").Value = 2296
$ws.Range("D27").Value = 173085
$ws.Range("E27").Value = 64609

# --- Ucrania (row 28) ------------------------------------------------------
$ws.Range("B28").Value = 208959
$ws.Range("C28").Value = 4027
$ws.Range("D28").Value = 92360
$ws.Range("E28").Value = 112470
$ws.Range("G28").Value = 64
$ws.Range("H28").Value = 4129

# --- Armenia (row 64) ------------------------------------------------------
$ws.Range("B64").Value = 50359
$ws.Range("C64").Value = 458
$ws.Range("D64").Value = 44001
$ws.Range("E64").Value = 5399
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 959

# --- Afganistan (row 71) ------------------------------------------------------
$ws.Range("B71").Value = 39268
$ws.Range("C71").Value = 14
$ws.Range("D71").Value = 32789
$ws.Range("E71").Value = 5021

# --- Hungria (row 81) ------------------------------------------------------
$ws.Range("B81").Value = 26461
$ws.Range("C81").Value = 894
$ws.Range("D81").Value = 5890
$ws.Range("E81").Value = 19806
$ws.Range("G81").Value = 8
$ws.Range("H81").Value = 765

# --- Georgia / Cabo Verde swap rank (rows 116-117) --------------------
# Georgia overtakes Cabo Verde (new case data) and moves to row 116,
# Cabo Verde (unchanged totals) drops to row 117.
$ws.Range("A116").Value = "Georgia"
$ws.Range("B116").Value = 6192
$ws.Range("C116").Value = 326
$ws.Range("D116").Value = 3120
$ws.Range("E116").Value = 3035
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 37

$ws.Range("A117").Value = "Cabo Verde"
$ws.Range("B117").Value = 5900
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 5228
$ws.Range("E117").Value = 613
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 59

# --- Taiwan (row 175) ------------------------------------------------------
$ws.Range("B175").Value = 514
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 483

# --- Islas Malvinas / Montserrat swap rank (rows 215-216) -------------
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
